$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date/time number format from A3 for new timestamp cells
$tsFormat = $ws.Cells.Item(3, 1).NumberFormat

# Row 3
$ws.Cells.Item(3, 1).Value = 45958.49250962963
$ws.Cells.Item(3, 1).NumberFormat = $tsFormat

# Row 4
$ws.Cells.Item(4, 1).Value = 45960.44648920139
$ws.Cells.Item(4, 2).Value = '-'
$ws.Cells.Item(4, 3).Value = 'add'
$ws.Cells.Item(4, 4).Value = 145
$ws.Cells.Item(4, 5).Value = 'テスト酒'
$ws.Cells.Item(4, 7).Value = '{}'
$ws.Cells.Item(4, 8).Value = '{''id'': 145, ''会員氏名'': ''テストくん'', ''name'': ''テスト酒'', ''蔵元'': '''', ''地域'': '''', ''category'': '''', ''精米歩合'': ''50'', ''updated_at'': Timestamp(''2025-10-30 10:42:56.512101''), ''備考'': ''''}'
$ws.Cells.Item(4, 1).NumberFormat = $tsFormat

# Row 5
$ws.Cells.Item(5, 1).Value = 45960.58553192129
$ws.Cells.Item(5, 2).Value = 'admin'
$ws.Cells.Item(5, 3).Value = 'add'
$ws.Cells.Item(5, 4).Value = 146
$ws.Cells.Item(5, 5).Value = 'テスト２'
$ws.Cells.Item(5, 7).Value = '{}'
$ws.Cells.Item(5, 8).Value = '{''id'': 146, ''会員氏名'': ''テストくん'', ''name'': ''テスト２'', ''蔵元'': '''', ''地域'': '''', ''category'': '''', ''精米歩合'': '''', ''updated_at'': Timestamp(''2025-10-30 14:03:09.842204''), ''備考'': '''', ''例会'': ''登録承認待ち''}'
$ws.Cells.Item(5, 1).NumberFormat = $tsFormat

# Row 6
$ws.Cells.Item(6, 1).Value = 45960.58623380787
$ws.Cells.Item(6, 2).Value = 'admin'
$ws.Cells.Item(6, 3).Value = 'update_meeting'
$ws.Cells.Item(6, 4).Value = 146
$ws.Cells.Item(6, 7).Value = '{}'
$ws.Cells.Item(6, 8).Value = '{''id'': 146, ''例会'': ''第24回''}'
$ws.Cells.Item(6, 1).NumberFormat = $tsFormat

# Row 7
$ws.Cells.Item(7, 1).Value = 45960.59161622685
$ws.Cells.Item(7, 2).Value = 'admin'
$ws.Cells.Item(7, 3).Value = 'update_meeting'
$ws.Cells.Item(7, 4).Value = 145
$ws.Cells.Item(7, 5).Value = 'テスト酒'
$ws.Cells.Item(7, 6).Value = '例会'
$ws.Cells.Item(7, 7).Value = '{''id'': 145, ''name'': ''テスト酒'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 10:42:56.512000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': ''50'', ''備考'': nan, ''例会'': nan, ''例会日時'': NaT}'
$ws.Cells.Item(7, 8).Value = '{''id'': 145, ''name'': ''テスト酒'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 10:42:56.512000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': ''50'', ''備考'': nan, ''例会'': ''24'', ''例会日時'': NaT}'
$ws.Cells.Item(7, 1).NumberFormat = $tsFormat

# Row 8
$ws.Cells.Item(8, 1).Value = 45960.5966872338
$ws.Cells.Item(8, 2).Value = 'admin'
$ws.Cells.Item(8, 3).Value = 'update_meeting'
$ws.Cells.Item(8, 4).Value = 145
$ws.Cells.Item(8, 5).Value = 'テスト酒'
$ws.Cells.Item(8, 6).Value = '例会'
$ws.Cells.Item(8, 7).Value = '{''id'': 145, ''name'': ''テスト酒'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 10:42:56.512000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': ''50'', ''備考'': nan, ''例会'': ''24'', ''例会日時'': NaT}'
$ws.Cells.Item(8, 8).Value = '{''id'': 145, ''name'': ''テスト酒'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 10:42:56.512000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': ''50'', ''備考'': nan, ''例会'': ''25'', ''例会日時'': NaT}'
$ws.Cells.Item(8, 1).NumberFormat = $tsFormat

# Row 9
$ws.Cells.Item(9, 1).Value = 45960.59791274305
$ws.Cells.Item(9, 2).Value = 'admin'
$ws.Cells.Item(9, 3).Value = 'update_meeting'
$ws.Cells.Item(9, 4).Value = 146
$ws.Cells.Item(9, 5).Value = 'テスト２'
$ws.Cells.Item(9, 6).Value = '例会'
$ws.Cells.Item(9, 7).Value = '{''id'': 146, ''name'': ''テスト２'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 14:03:09.842000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': nan, ''備考'': nan, ''例会'': ''第24回'', ''例会日時'': NaT}'
$ws.Cells.Item(9, 8).Value = '{''id'': 146, ''name'': ''テスト２'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 14:03:09.842000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': nan, ''備考'': nan, ''例会'': ''26'', ''例会日時'': NaT}'
$ws.Cells.Item(9, 1).NumberFormat = $tsFormat

# Row 10
$ws.Cells.Item(10, 1).Value = 45960.60267283207
$ws.Cells.Item(10, 2).Value = 'admin'
$ws.Cells.Item(10, 3).Value = 'delete'
$ws.Cells.Item(10, 4).Value = 146
$ws.Cells.Item(10, 5).Value = 'テスト２'
$ws.Cells.Item(10, 7).Value = '{''id'': 146, ''name'': ''テスト２'', ''category'': nan, ''quantity'': nan, ''updated_at'': Timestamp(''2025-10-30 14:03:09.842000''), ''会員氏名'': ''テストくん'', ''蔵元'': nan, ''地域'': nan, ''精米歩合'': nan, ''備考'': nan, ''例会'': 26, ''例会日時'': NaT}'
$ws.Cells.Item(10, 8).Value = '{}'
$ws.Cells.Item(10, 1).NumberFormat = $tsFormat
